$d = $word.ActiveDocument

# --- Edit 1: paragraph "1. Show a documentation header. " ---
# Split the trailing space off, then add the new sentence plus a trailing space run.
$d.Content.Find.Execute(
    "1. Show a documentation header. ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "1. Show a documentation header.", 2)

# Locate that paragraph again and append the two new runs after its current text.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq "1. Show a documentation header.`r") {
        $p1 = $cand
        break
    }
}
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.InsertAfter(" Include your name, the purpose of the program, and the conditions under which others may or may not use your results.")
$r1.InsertAfter(" ")

# --- Edit 2: paragraph "2. D" gets the rest of the sentence, then nine new Q&A blocks ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq "2. D`r") {
        $p2 = $cand
        break
    }
}
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)

$newBlock = "ownload and import the file samara.txt. Review the data dictionary, if needed. Display the first ten rows of data below." `
    + "`r`r3. There appears to be a relationship between load and velocity. Draw a graph illustrating this relationship with load on the X axis and velocity on the Y axis. Make sure that your graph follows good design principles. Include a trendline in your graph. Describe the graph and highlight unusual values (if any)." `
    + "`r`r4. Calculate a simple linear regression (Analyze | Regression | Linear) using load as the independent variable and velocity as the dependent variable. Show the table with R squared and interpret this number." `
    + "`r`r5. Show the analysis of variance table and use the F ratio to test the hypothesis that the population slope is zero. Interpret your results." `
    + "`r`r6. Compute the residuals and examine the normality assumption using a histogram. Interpret your results." `
    + "`r`r7. The researchers believe that the average velocity is the same for all three trees in the study. Display a boxplot that would provide an informal assessment of this hypothesis." `
    + "`r`r8. Does the boxplot show evidence of unequal variances or non-normality?" `
    + "`r`r9. Show the analysis of variance table for examining the average velocity across all three trees. Explain what the F ratio is testing and interpret your results." `
    + "`r`r10. Since the sample sizes in each group are almost equal, you can use the Tukey post hoc test to examine which trees differ from the others. Display the SPSS output associated with this test and interpret the results."

$r2.InsertAfter($newBlock)

# --- Split run 3 into two runs (graph sentence / describe sentence) ---
$d.Content.Find.Execute(
    "Include a trendline in your graph. Describe the graph", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Include a trendline in your graphZZZSPLITZZZ. Describe the graph", 2)
$d.Content.Find.Execute("ZZZSPLITZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Split run 10 into two runs (intro sentence / closing sentence) ---
$d.Content.Find.Execute(
    "associated with this test and interpret", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "associated with ZZZSPLITZZZthis test and interpret", 2)
$d.Content.Find.Execute("ZZZSPLITZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
